# Add 14 new rows of "Sword World 2.0" collectibles data to the checklist,
# matching the column-by-column order the data was originally typed in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 - Players Handbook Faydan Museum
$ws.Cells.Item(24, 1).Value = 2012
$ws.Cells.Item(24, 5).Value = "faydan_museum.jpg"
$ws.Cells.Item(24, 2).Value = "プレイヤーズ・ハンドブック フェイダン博物誌"
$ws.Cells.Item(24, 3).Value = "Players Handbook Faydan Museum"
$ws.Cells.Item(24, 4).Value = "Fujimi Shobo"
$ws.Cells.Item(24, 6).Value = "supplement"

# Row 25 - Players Handbook Salz Museum
$ws.Cells.Item(25, 1).Value = 2012
$ws.Cells.Item(25, 5).Value = "salz_museum.jpg"
$ws.Cells.Item(25, 2).Value = "プレイヤーズ・ハンドブック ザルツ博物誌"
$ws.Cells.Item(25, 3).Value = "Players Handbook Salz Museum"
$ws.Cells.Item(25, 4).Value = "Kadokawa"
$ws.Cells.Item(25, 6).Value = "supplement"

# Row 26 - Players Handbook Eurelia Natural History
$ws.Cells.Item(26, 1).Value = 2013
$ws.Cells.Item(26, 2).Value = "プレイヤーズ・ハンドブック ユーレリア博物誌"
$ws.Cells.Item(26, 3).Value = "Players Handbook Eurelia Natural History"
$ws.Cells.Item(26, 4).Value = "Kadokawa"
$ws.Cells.Item(26, 5).Value = "eurelia.jpg"
$ws.Cells.Item(26, 6).Value = "supplement"

# Row 27 - Players Handbook Dagnia Museum
$ws.Cells.Item(27, 1).Value = 2015
$ws.Cells.Item(27, 2).Value = "プレイヤーズ・ハンドブック ダグニア博物誌"
$ws.Cells.Item(27, 3).Value = "Players Handbook Dagnia Museum"
$ws.Cells.Item(27, 4).Value = "Kadokawa"
$ws.Cells.Item(27, 5).Value = "dagnia_museum.jpg"
$ws.Cells.Item(27, 6).Value = "supplement"

# Row 28 - Players Handbook Dilfram Museum
$ws.Cells.Item(28, 1).Value = 2016
$ws.Cells.Item(28, 2).Value = "プレイヤーズ・ハンドブック ディルフラム博物誌 "
$ws.Cells.Item(28, 3).Value = "Players Handbook Dilfram Museum"
$ws.Cells.Item(28, 4).Value = "Kadokawa"
$ws.Cells.Item(28, 5).Value = "dilfram_museum.jpg"
$ws.Cells.Item(28, 6).Value = "supplement"

# Row 29 - Sword World 2.0 Tour 1: Luferia
$ws.Cells.Item(29, 1).Value = 2009
$ws.Cells.Item(29, 2).Value = "ソード・ワールド2.0 ツアー(1) ルーフェリア"
$ws.Cells.Item(29, 3).Value = "Sword World 2.0 Tour 1: Luferia"
$ws.Cells.Item(29, 4).Value = "Fujimi Shobo"
$ws.Cells.Item(29, 5).Value = "luferia.jpg"
$ws.Cells.Item(29, 6).Value = "supplement"

# Row 30 - Sword World 2.0 Tour 2: Lios
$ws.Cells.Item(30, 1).Value = 2010
$ws.Cells.Item(30, 2).Value = "ソード・ワールド2.0ツアー (2)　リオス"
$ws.Cells.Item(30, 3).Value = "Sword World 2.0 Tour 2: Lios"
$ws.Cells.Item(30, 4).Value = "Fujimi Shobo"
$ws.Cells.Item(30, 5).Value = "lios.jpg"
$ws.Cells.Item(30, 6).Value = "supplement"

# Row 31 - Battle Campaign Book Calzoral's Magic Angel
$ws.Cells.Item(31, 1).Value = 2014
$ws.Cells.Item(31, 5).Value = "calzorals_magic_angel.jpg"
$ws.Cells.Item(31, 2).Value = "バトルキャンペーンブックカルゾラルの魔動天使"
$ws.Cells.Item(31, 3).Value = "Battle Campaign Book Calzoral's Magic Angel"
$ws.Cells.Item(31, 4).Value = "Kadokawa"
$ws.Cells.Item(31, 6).Value = "supplement"

# Row 32 - Battle Campaign Book Procercia Secret History-Dawn Princess
$ws.Cells.Item(32, 1).Value = 2015
$ws.Cells.Item(32, 2).Value = "バトルキャンペーンブック プロセルシア秘史 ―暁をうたう竜の姫"
$ws.Cells.Item(32, 3).Value = "Battle Campaign Book Procercia Secret History-Dawn Princess"
$ws.Cells.Item(32, 4).Value = "Kadokawa"
$ws.Cells.Item(32, 5).Value = "procercia_secret_history.jpg"
$ws.Cells.Item(32, 6).Value = "supplement"

# Row 33 - Story & Data Book Dragon Raid Senryuden
$ws.Cells.Item(33, 1).Value = 2014
$ws.Cells.Item(33, 2).Value = "ストーリー&データブックドラゴンレイド戦竜伝 "
$ws.Cells.Item(33, 3).Value = "Story & Data Book Dragon Raid Senryuden"
$ws.Cells.Item(33, 4).Value = "Kadokawa"
$ws.Cells.Item(33, 5).Value = "dragon_raid_senryuden.jpg"
$ws.Cells.Item(33, 6).Value = "supplement"

# Row 34 - Grand Campaign Dragon Raid Begins-White Dragon Maiden-
$ws.Cells.Item(34, 1).Value = 2014
$ws.Cells.Item(34, 2).Value = "グランドキャンペーン ドラゴンレイドビギンズ ‐白き竜の乙女‐"
$ws.Cells.Item(34, 3).Value = "Grand Campaign Dragon Raid Begins-White Dragon Maiden-"
$ws.Cells.Item(34, 4).Value = "Kadokawa"
$ws.Cells.Item(34, 5).Value = "white_dragon_maiden.jpg"
$ws.Cells.Item(34, 6).Value = "supplement"

# Row 35 - Story & Data Book Dragon Raid Senryuden (2)
$ws.Cells.Item(35, 1).Value = 2015
$ws.Cells.Item(35, 2).Value = "ストーリー&データブック ドラゴンレイド戦竜伝 (2)"
$ws.Cells.Item(35, 3).Value = "Story & Data Book Dragon Raid Senryuden (2)"
$ws.Cells.Item(35, 4).Value = "Kadokawa"
$ws.Cells.Item(35, 5).Value = "dragon_raid_senryuden2.jpg"
$ws.Cells.Item(35, 6).Value = "supplement"

# Row 36 - Scenario Collection: Challenge! Labyrinth Called by the Magic Sword
# (image filename for this row was filled in later, after row 37 below)
$ws.Cells.Item(36, 1).Value = 2008
$ws.Cells.Item(36, 2).Value = "シナリオ集　挑戦! 魔剣が呼ぶ迷宮"
$ws.Cells.Item(36, 3).Value = "Scenario Collection: Challenge! Labyrinth Called by the Magic Sword"
$ws.Cells.Item(36, 4).Value = "Fujimi Shobo"
$ws.Cells.Item(36, 6).Value = "scenario"

# Row 37 - Scenario Collection (2) Fengyun! A City Where the Singing Voice Echoes
$ws.Cells.Item(37, 1).Value = 2008
$ws.Cells.Item(37, 2).Value = "シナリオ集(2) 風雲!歌声が響く都市"
$ws.Cells.Item(37, 3).Value = "Scenario Collection (2) Fengyun! A City Where the Singing Voice Echoes"
$ws.Cells.Item(37, 4).Value = "Fujimi Shobo"
$ws.Cells.Item(37, 5).Value = "fengyun.jpg"
$ws.Cells.Item(37, 6).Value = "scenario"

# Back-fill the row 36 image now that it's known
$ws.Cells.Item(36, 5).Value = "labyrinth_magic_sword.jpg"

$ws.Range("F37").Select()
